# Update scripts with new TPM: recomputed NATMI ligand/receptor
# expression, specificity, and edge-weight statistics for the
# Cdh1 (MuSCs / Resolving-Mac) -> Igf1r (ECs / FAPs / MuSCs /
# Resolving-Mac) pairs on rows 2-9, columns E-J and M-T.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1443736666666667
$ws.Range("H2").Value = 0.433121
$ws.Range("I2").Value = 0.7378778224885942
$ws.Range("J2").Value = 0.7378778224885942
$ws.Range("M2").Value = 15.03663066666667
$ws.Range("N2").Value = 45.109892
$ws.Range("O2").Value = 0.279146411176606
$ws.Range("P2").Value = 0.279146411176606
$ws.Range("Q2").Value = 2.170893503659111
$ws.Range("R2").Value = 19.538041532932
$ws.Range("S2").Value = 0.2059759460344999
$ws.Range("T2").Value = 0.2059759460344998

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1443736666666667
$ws.Range("H3").Value = 0.433121
$ws.Range("I3").Value = 0.7378778224885942
$ws.Range("J3").Value = 0.7378778224885942
$ws.Range("O3").Value = 0.2673306493381863
$ws.Range("P3").Value = 0.2673306493381863
$ws.Range("Q3").Value = 2.079003514790223
$ws.Range("R3").Value = 18.711031633112
$ws.Range("S3").Value = 0.1972573574181229
$ws.Range("T3").Value = 0.1972573574181229

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1443736666666667
$ws.Range("H4").Value = 0.433121
$ws.Range("I4").Value = 0.7378778224885942
$ws.Range("J4").Value = 0.7378778224885942
$ws.Range("M4").Value = 22.16851266666667
$ws.Range("N4").Value = 66.505538
$ws.Range("O4").Value = 0.411545703901694
$ws.Range("P4").Value = 0.411545703901694
$ws.Range("Q4").Value = 3.200549458233112
$ws.Range("R4").Value = 28.804945124098
$ws.Range("S4").Value = 0.3036704478495177
$ws.Range("T4").Value = 0.3036704478495177

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1443736666666667
$ws.Range("H5").Value = 0.433121
$ws.Range("I5").Value = 0.7378778224885942
$ws.Range("J5").Value = 0.7378778224885942
$ws.Range("M5").Value = 2.261165333333333
$ws.Range("N5").Value = 6.783496
$ws.Range("O5").Value = 0.04197723558351375
$ws.Range("P5").Value = 0.04197723558351374
$ws.Range("Q5").Value = 0.3264527301128889
$ws.Range("R5").Value = 2.938074571016
$ws.Range("S5").Value = 0.03097407118645386
$ws.Range("T5").Value = 0.03097407118645385

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.051287
$ws.Range("H6").Value = 0.153861
$ws.Range("I6").Value = 0.2621221775114058
$ws.Range("J6").Value = 0.2621221775114058
$ws.Range("M6").Value = 15.03663066666667
$ws.Range("N6").Value = 45.109892
$ws.Range("O6").Value = 0.279146411176606
$ws.Range("P6").Value = 0.279146411176606
$ws.Range("Q6").Value = 0.7711836770013334
$ws.Range("R6").Value = 6.940653093012
$ws.Range("S6").Value = 0.07317046514210619
$ws.Range("T6").Value = 0.07317046514210619

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.051287
$ws.Range("H7").Value = 0.153861
$ws.Range("I7").Value = 0.2621221775114058
$ws.Range("J7").Value = 0.2621221775114058
$ws.Range("O7").Value = 0.2673306493381863
$ws.Range("P7").Value = 0.2673306493381863
$ws.Range("Q7").Value = 0.7385408691546667
$ws.Range("R7").Value = 6.646867822392001
$ws.Range("S7").Value = 0.07007329192006344
$ws.Range("T7").Value = 0.07007329192006345

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.051287
$ws.Range("H8").Value = 0.153861
$ws.Range("I8").Value = 0.2621221775114058
$ws.Range("J8").Value = 0.2621221775114058
$ws.Range("M8").Value = 22.16851266666667
$ws.Range("N8").Value = 66.505538
$ws.Range("O8").Value = 0.411545703901694
$ws.Range("P8").Value = 0.411545703901694
$ws.Range("Q8").Value = 1.136956509135333
$ws.Range("R8").Value = 10.232608582218
$ws.Range("S8").Value = 0.1078752560521763
$ws.Range("T8").Value = 0.1078752560521763

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.051287
$ws.Range("H9").Value = 0.153861
$ws.Range("I9").Value = 0.2621221775114058
$ws.Range("J9").Value = 0.2621221775114058
$ws.Range("M9").Value = 2.261165333333333
$ws.Range("N9").Value = 6.783496
$ws.Range("O9").Value = 0.04197723558351375
$ws.Range("P9").Value = 0.04197723558351374
$ws.Range("Q9").Value = 0.1159683864506666
$ws.Range("R9").Value = 1.043715478056
$ws.Range("S9").Value = 0.01100316439705989
$ws.Range("T9").Value = 0.01100316439705989
